$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P2").Value = 1.173940532471115
$ws.Range("P3").Value = 0.6687495246397893
$ws.Range("P4").Value = 0.6450408693787498
$ws.Range("P5").Value = 0.6581965283700736
$ws.Range("P6").Value = 0.7599649704893262
$ws.Range("P7").Value = 1.526723262206673
$ws.Range("P8").Value = 0.6571205391460256
$ws.Range("P9").Value = 0.6334061194151843
$ws.Range("P10").Value = 0.7033441035257412
$ws.Range("P11").Value = 0.7651462897421865
$ws.Range("P12").Value = 0.7580490425937152
$ws.Range("P13").Value = 0.8571530970826886
$ws.Range("P14").Value = 0.8172362399217724
$ws.Range("P15").Value = 0.5838850344709698
$ws.Range("P16").Value = 0.6350345065702816
$ws.Range("P17").Value = 0.6002439556025914
$ws.Range("P18").Value = 0.5742170078473631
$ws.Range("P19").Value = 0.7177770347464668
$ws.Range("P20").Value = 1.205769495845655
$ws.Range("P21").Value = 0.5594446828765947
$ws.Range("P22").Value = 0.6052881758016058
$ws.Range("P23").Value = 0.6612256998766306
$ws.Range("P24").Value = 0.6207933488881647
$ws.Range("P25").Value = 0.6709494555951921
$ws.Range("P26").Value = 0.5945000833547658
$ws.Range("P27").Value = 0.6690940401230234
$ws.Range("P28").Value = 1.386192599252586
$ws.Range("P29").Value = 3.209388084453111
$ws.Range("P30").Value = 2.875728074034313
$ws.Range("P31").Value = 2.536223211262842
$ws.Range("P32").Value = 3.209680883537202
$ws.Range("P33").Value = 2.951980627688052
$ws.Range("P34").Value = 2.581096709220017
$ws.Range("P35").Value = 1.944977867088071
$ws.Range("P36").Value = 1.902997215914106
$ws.Range("P37").Value = 0.9311757893288353
$ws.Range("P38").Value = 0.7639131017440088
$ws.Range("P39").Value = 0.7865765523628162
$ws.Range("P40").Value = 2.08999422632383
$ws.Range("P41").Value = 0.6512220892168789
$ws.Range("P42").Value = 0.7214535434586795
$ws.Range("P43").Value = 0.5752907435528255
$ws.Range("P44").Value = 1.088286471391299
$ws.Range("P45").Value = 0.5912974701318487
$ws.Range("P46").Value = 0.5819542572851222
$ws.Range("P47").Value = 0.6378035955409712
$ws.Range("P48").Value = 0.5757688880479208
$ws.Range("P49").Value = 0.591550316930182
$ws.Range("P50").Value = 0.5679457863888585
$ws.Range("P51").Value = 0.9753587577666069
$ws.Range("P52").Value = 1.223267635565539
$ws.Range("P53").Value = 0.7317666296063944
